# Generate Report for Handback
# Update the handoff/handback timestamps on the per-language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 08:32:47"
$wsZhCn.Range("H2").Value = "2016-03-19 08:33:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 08:32:55"
$wsDeDe.Range("H2").Value = "2016-03-19 08:34:02"
